# Big stimulus update:
#  - normalize "correct_ans" (column L) codes to full words (r -> right, b -> center, y -> left)
#  - swap "face" image stimuli for "book" image stimuli wherever they occur

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count
$lastCol = $ur.Columns.Count

# 1) Expand the abbreviated answer codes in column L (correct_ans) to full words,
#    for every data row (row 1 is the header "correct_ans").
$ansMap = @{ "r" = "right"; "b" = "center"; "y" = "left" }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 12)
    $val = $cell.Value()
    if (($val -is [string]) -and $ansMap.ContainsKey($val)) {
        $cell.Value = $ansMap[$val]
    }
}

# 2) Replace every "face//face_NN.jpg" stimulus filename with "book//book_NN.jpg"
#    wherever it appears across the whole used range.
for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if (($val -is [string]) -and $val.Contains("face//face_")) {
            $cell.Value = $val.Replace("face//face_", "book//book_")
        }
    }
}
